# Generate Report for Handback
# Refreshes the handback-status report timestamps produced by a new
# handback-generation run (ae028673-* row picks up the new run's
# handoff/handback datetimes; cd54ed06-* row is untouched by this run).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for the
# ae028673 row (row 2) advances to the new generation timestamp.
$overview.Range("G2").Value = "2016-08-28 06:53:06"

# zh-cn sheet: row 2 (ae028673) Correspond Handoff / Handback datetimes
# move forward to reflect the new handback run.
$zhcn.Range("H2").Value = "2016-08-28 06:52:58"
$zhcn.Range("K2").Value = "2016-08-28 06:53:26"

# de-de sheet: row 2 (ae028673) Correspond Handoff / Handback datetimes
# move forward to reflect the new handback run.
$dede.Range("H2").Value = "2016-08-28 06:53:06"
$dede.Range("K2").Value = "2016-08-28 06:53:32"
